$d = $word.ActiveDocument

# "...too low.." (the sentence's trailing period plus a stray extra period run)
# becomes "...too low…" (drop one period, turn the remaining one into an ellipsis).
$d.Content.Find.Execute("too low..", $true, $false, $false, $false, $false, `
                         $true, 1, $false, "too low…", 2)
